# Tripadvisor New Orleans shard 137 update:
#  - swap the two sheets' identity/order: review_info becomes sheet #1
#    (sheetId 1), hotel_info becomes sheet #2 (sheetId 2)
#  - review_info keeps its original 25 columns (A:Y), no data rows
#  - hotel_info gains a new "State" column (right after Hotel_Name) and
#    its single data row gets "Louisiana" in that new column

$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item(1)   # currently "hotel_info"
$wsB = $wb.Worksheets.Item(2)   # currently "review_info"

# wipe all existing content/shared-string references on both sheets first
$wsA.Cells.Clear()
$wsB.Cells.Clear()

# swap the sheet names (go through a temp name to dodge the collision)
$wsA.Name = "__tmp_sheet_swap__"
$wsB.Name = "hotel_info"
$wsA.Name = "review_info"

$wsReview = $wsA   # now named "review_info", sheetId 1 / rId1 / physically sheet1.xml
$wsHotel  = $wsB   # now named "hotel_info",  sheetId 2 / rId2 / physically sheet2.xml

# ---- review_info header row (A1:Y1) ----
$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)
for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $wsReview.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# ---- hotel_info header row (A1:J1) ----
$hotelHeaders = @(
    "STR",
    "Hotel_Name",
    "State",
    "City",
    "Zip",
    "TA_ReviewURL",
    "Tripadvisor_Hotel_Name",
    "English_Reviews_num",
    "Local_Rank",
    "Total_Reviews_num"
)
for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $wsHotel.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

# ---- hotel_info data row (A2:J2) ----
$wsHotel.Cells.Item(2, 1).Value = 15984
$wsHotel.Cells.Item(2, 2).Value = "Le Pavillon Hotel"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"
$wsHotel.Cells.Item(2, 4).Value = "New Orleans"
$wsHotel.Cells.Item(2, 5).Value = 70112
$wsHotel.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g60864-d89091-Reviews-Le_Pavillon_Hotel-New_Orleans_Louisiana.html"
$wsHotel.Cells.Item(2, 7).Value = "Le Pavillon Hotel"

# these three look numeric but must stay text, like in the original file
$wsHotel.Cells.Item(2, 8).NumberFormat = "@"
$wsHotel.Cells.Item(2, 8).Value = "2896"
$wsHotel.Cells.Item(2, 9).NumberFormat = "@"
$wsHotel.Cells.Item(2, 9).Value = "106"
$wsHotel.Cells.Item(2, 10).NumberFormat = "@"
$wsHotel.Cells.Item(2, 10).Value = "2976"
